# Revert commit: undo the merge that had (a) appended a "?" question mark
# removal on the "Questões" sheet label, (b) checked off answers in the
# Questões grid, and (c) left the "Questões" tab active.
#
# This script restores the previous state: adds back the missing "?" to the
# question text, clears every answer mark in the Questões grid (so the
# dependent Ponto de Função calculations fall back to their unanswered
# values), and returns the active tab / selections to the "Resultados" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet references (by fixed tab position, robust to accented names) ---
$wsResultados = $wb.Worksheets.Item(1)   # "Resultados"
$wsCalculo    = $wb.Worksheets.Item(2)   # "Cálculo do Ponto de Função"
$wsQuestoes   = $wb.Worksheets.Item(4)   # "Questões"

# --- 1. Fix the question text: append the missing "?" ---
$cellB10 = $wsQuestoes.Range("B10")
$cellB10.Value = $cellB10.Value() + "?"

# --- 2. Clear every selected answer mark in the Questões grid (E3:J16) ---
$answerCells = @("J3", "E4", "E5", "H6", "G7", "E8", "E9", "E10", "F11", "F12", "G13", "J14", "J15", "H16")
foreach ($addr in $answerCells) {
    $wsQuestoes.Range($addr).ClearContents()
}

# --- 3. Restore view state: selections + active tab back on "Resultados" ---
$wsCalculo.Range("D17").Select()
$wsQuestoes.Range("B11:D11").Select()

$wsResultados.Activate()
$wsResultados.Range("E22").Select()
